$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 6 ("SOURCES OF FINANCE") switches to a different built-in
#    table style (tableStyleId GUID changes).
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$tableShape = $slide6.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{8C4483A6-5E51-4526-AC9F-D332D6540FA7}")

# ---------------------------------------------------------------------------
# 2) The deck's theme switches from the "Integral" design back to the
#    default "Office Theme" palette. Push the Office Theme's 12 scheme
#    colors (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) onto the
#    presentation's live theme color scheme.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

$officeThemeRGB = @(
    0x000000,  # 1  dk1
    0xFFFFFF,  # 2  lt1
    0x6A5444,  # 3  dk2
    0xE6E6E7,  # 4  lt2
    0xD59B5B,  # 5  accent1
    0x317DED,  # 6  accent2
    0xA5A5A5,  # 7  accent3
    0x00C0FF,  # 8  accent4
    0xC47244,  # 9  accent5
    0x47AD70,  # 10 accent6
    0xC16305,  # 11 hlink
    0x724F95   # 12 folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeRGB[$i - 1]
}
